# Update the build timestamp embedded in the workbook text from
# "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---
$cellA2 = $aboutSheet.Range("A2")
$cellA2.Value = $cellA2.Text.Replace($oldStamp, $newStamp)

$cellA6 = $aboutSheet.Range("A6")
$cellA6.Value = $cellA6.Text.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet (column S, rows 2-11) ---
for ($row = 2; $row -le 11; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # Column S = 19
    $cell.Value = $cell.Text.Replace($oldStamp, $newStamp)
}
